$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.916.90"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.158.92"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'246.35"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").Value = "'65.01"
$ws.Range("E7").Value = "  -8.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("E9").Value = "  -4.07%  "
$ws.Range("D10").Value = "'59.21"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").Value = "'34.89"
$ws.Range("E12").Value = "  -16.14%  "
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").Value = "'6.76"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").Value = "2.474.53"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'14.11"
$ws.Range("E16").Value = "  -5.78%  "
$ws.Range("D17").Value = "'0.839"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "2.173.69"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "40.815.86"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").Value = "0.0₃0928"
$ws.Range("E20").Value = "  -4.45%  "
$ws.Range("D21").Value = "'70.93"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("D23").Value = "'228.00"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("E24").Value = "  -6.53%  "
$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D25").Value = "'3.81"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'11.09"
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("D28").Value = "'2.39"
$ws.Range("E28").Value = "  -5.85%  "
$ws.Range("D29").Value = "'3.71"
$ws.Range("E29").Value = "  -5.89%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'166.89"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'20.03"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "'5.59"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'0.0737"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("D38").Value = "'3.92"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'24.20"
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("D40").Value = "'0.0297"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").Value = "'5.40"
$ws.Range("E42").Value = "  -9.35%  "
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'59.82"
$ws.Range("E44").Value = "  -13.22%  "
$ws.Range("D45").Value = "'11.02"
$ws.Range("E45").Value = "  -7.29%  "
$ws.Range("E46").Value = "  -9.84%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'8.37"
$ws.Range("E48").Value = "  -5.06%  "
$ws.Range("D49").Value = "'0.0979"
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").Value = "'1.13"
$ws.Range("E51").Value = "  -5.70%  "
